# Update "想去人数" (want-to-go count) values in column F on the "展览"
# and "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Row -> new F value mapping shared by both worksheets (same events, same
# relative order, but different absolute row numbers on sheet 全部类型
# because it interleaves rows from other sheets).
$updatesSheet1 = @{
    5  = 15901
    8  = 713
    11 = 9104
    12 = 396
    16 = 207
    20 = 66
    21 = 571
    24 = 65
    25 = 1120
    27 = 17
    28 = 25
    29 = 498
    33 = 66
    36 = 332
    39 = 5597
}

$updatesSheet4 = @{
    5  = 15901
    8  = 713
    11 = 9104
    12 = 396
    16 = 207
    20 = 66
    21 = 571
    24 = 65
    25 = 1120
    27 = 17
    28 = 25
    29 = 498
    35 = 66
    38 = 332
    41 = 5597
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Range("F$row").Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Range("F$row").Value = $updatesSheet4[$row]
}
